$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9546350836753845
$ws.Range("B1").Value = 2.011213302612305
$ws.Range("C1").Value = 4.127926349639893
$ws.Range("D1").Value = 3.218059301376343
$ws.Range("E1").Value = 1.430335164070129
